$d = $word.ActiveDocument

$p1Text = 'המאמר היומי של מייק ואוראל: 19.01.25' + [char]11 + 'The Lottery Ticket Hypothesis: Finding Sparse, Trainable Neural Networks'
# Paragraph 1: date/author line + paper title (two runs separated by a line break)
$rng = ($d.Paragraphs(1)).Range
$textRange = $d.Range($rng.Start, $rng.End - 1)
$textRange.Select()
if ($textRange.Start -lt $textRange.End) {
    $word.Selection.Delete()
}
$word.Selection.TypeText($p1Text)

# Paragraph 2: intro of Lottery Ticket Hypothesis
$rng = ($d.Paragraphs(2)).Range
$textRange = $d.Range($rng.Start, $rng.End - 1)
$textRange.Select()
if ($textRange.Start -lt $textRange.End) {
    $word.Selection.Delete()
}
$word.Selection.TypeText('היפותזת כרטיס הלוטו (Lottery Ticket Hypothesis) אומרת שבתוך רשת נוירונים  צפופה (dense neural nets) המאותחלת בצורה רנדומלית, יש תת-רשת (או "כרטיס מנצח") שמאמנים אותה בנפרד, היא יכולה להגיע לביצועים כמו של הרשת המקורית.')

# Paragraph 3: pruning finds sub-networks
$rng = ($d.Paragraphs(3)).Range
$textRange = $d.Range($rng.Start, $rng.End - 1)
$textRange.Select()
if ($textRange.Start -lt $textRange.End) {
    $word.Selection.Delete()
}
$word.Selection.TypeText('נמצא שטכניקת חיתוך(pruning) סטנדרטית מגלה באופן טבעי תת-רשתות כאלה, אשר עבורן מתקיים כי האתחול המחודש תחת אותם hyperparameters, משמר את התוצאות של הרשת המקורית בעלות זולה יותר, כך שהכרטיסים המנצחים הם תת-רשתות אשר "זכו בהגרלת האתחול", ובהן המשקלים ההתחלתיים הופכים את האימון לאפקטיבי במיוחד.')

# Paragraph 4: importance of initial weights
$rng = ($d.Paragraphs(4)).Range
$textRange = $d.Range($rng.Start, $rng.End - 1)
$textRange.Select()
if ($textRange.Start -lt $textRange.End) {
    $word.Selection.Delete()
}
$word.Selection.TypeText('הרעיון הזה מדגיש את החשיבות של המשקלים ההתחלתיים של הרשת. הכרטיסים המנצחים אינם תת-רשתות אקראיות, אלא כאלה שמתאימות במיוחד בגלל האתחול שלהן. תהליך מציאת התת-רשתות הללו אינו פשוט, כיוון שהוא כרוך בזיהוי החלקים הקריטיים(הנוירונים המשמעותיים) ברשת כבר מההתחלה.')

# Paragraph 5: heading - what is network pruning
$rng = ($d.Paragraphs(5)).Range
$textRange = $d.Range($rng.Start, $rng.End - 1)
$textRange.Select()
if ($textRange.Start -lt $textRange.End) {
    $word.Selection.Delete()
}
$word.Selection.TypeText('מה זה חיתוך רשת?')

# Paragraph 6: definition of pruning
$rng = ($d.Paragraphs(6)).Range
$textRange = $d.Range($rng.Start, $rng.End - 1)
$textRange.Select()
if ($textRange.Start -lt $textRange.End) {
    $word.Selection.Delete()
}
$word.Selection.TypeText('חיתוך (Pruning) הוא טכניקה המסירה משקלים לא חשובים מרשת הנוירונים. לפי היפותזת כרטיס הלוטו, החיתוך עוזר לייעל את הרשת בכך שהוא מסיר נוירונים וחיבורים מיותרים, וכך יוצר רשת קלה, מהירה ויעילה יותר, ששומרת על הביצועים של הרשת המקורית ולעיתים אף משפרת אותם. החיתוך חושף את "הכרטיסים המנצחים": בתחילה, הרשת מכילה יותר מדי פרמטרים (רשת גדולה וצפופה), ואז במהלך האימון והחיתוך של המשקלים הלא משמעותיים, תת-הרשתות היעילות האלו מתגלות.')

$p7Text = 'סוגי חיתוך' + [char]11 + 'חיתוך לא מובנה (Unstructured Pruning): כאן אפשר להסיר כל משקל או קבוצה של משקלים, ללא מגבלות. זה יוצר רשת נוירונים "דלילה" שבה רק חלק מהמשקלים נשארים. טכניקה זו נקראת גם חיתוך משקלים (Weight Pruning). בחיתוך שכזה, אין בחירה מוגדרת מראש מה ייחתך, הכל לפי הבחירה הפחותה ביותר של התרומה של אותו נוירון שנבחר להיחתך.'
# Paragraph 7: "Types of pruning" heading + unstructured pruning text (line break)
$rng = ($d.Paragraphs(7)).Range
$textRange = $d.Range($rng.Start, $rng.End - 1)
$textRange.Select()
if ($textRange.Start -lt $textRange.End) {
    $word.Selection.Delete()
}
$word.Selection.TypeText($p7Text)

# Paragraph 8: structured-pruning paragraph replaces the old arxiv link
$rng = ($d.Paragraphs(8)).Range
$textRange = $d.Range($rng.Start, $rng.End - 1)
$textRange.Select()
if ($textRange.Start -lt $textRange.End) {
    $word.Selection.Delete()
}
$word.Selection.TypeText('חיתוך מובנה (Structured Pruning): כאן מסירים קבוצות שלמות של משקלים, כמו נוירונים שלמים ברשת קדמית (FFN). התוצאה היא רשת נוירונים "צפופה" אך קטנה יותר. הבחירה כאן היא מושכלת, בה המבניות של הרשת חשובה להישמר, יכול להיות שיהיה נוירון שלא יבחר להיחתך על מנת לא לפגוע במבניות שנבחרה, לעומת נוירונים אחרים.')

# Append new paragraphs: one-shot vs iterative heading(+break+text), iterative text, new arxiv link
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newCount = $d.Paragraphs.Count
$p9Text = 'חיתוך בבת אחת מול חיתוך איטרטיבי' + [char]11 + 'חיתוך בבת אחת (One-shot Pruning): מאמנים את הרשת פעם אחת, חותכים אחוז מסוים מהמשקלים (p%), ואז מאתחלים מחדש את המשקלים שנשארו. מדובר בהנחה כי באיטרציה אחת הגענו לפתרון הסופי והמיוחל, ללא צורך בתהליך חוזר ומתמשך.'
# New paragraph 9: one-shot vs iterative heading + one-shot text
$rng = ($d.Paragraphs($newCount)).Range
$textRange = $d.Range($rng.Start, $rng.End - 1)
$textRange.Select()
if ($textRange.Start -lt $textRange.End) {
    $word.Selection.Delete()
}
$word.Selection.TypeText($p9Text)

$p9 = $d.Paragraphs($newCount)
$p9.Range.InsertParagraphAfter()
$newCount2 = $d.Paragraphs.Count
# New paragraph 10: iterative pruning text
$rng = ($d.Paragraphs($newCount2)).Range
$textRange = $d.Range($rng.Start, $rng.End - 1)
$textRange.Select()
if ($textRange.Start -lt $textRange.End) {
    $word.Selection.Delete()
}
$word.Selection.TypeText('חיתוך איטרטיבי (Iterative Pruning): מאמנים את הרשת, חותכים חלק מהמשקלים, מאתחלים מחדש, וחוזרים על התהליך כמה פעמים. בכל סיבוב חותכים אחוז קטן מהמשקלים ששרדו מהסיבוב הקודם. תוצאות מראות שחיתוך איטרטיבי מצליח למצוא כרטיסים מנצחים שמגיעים לאותם ביצועים כמו של הרשת המקורית, תוך שימוש ברשת קטנה יותר בהשוואה לחיתוך בבת אחת.')

$p10 = $d.Paragraphs($newCount2)
$p10.Range.InsertParagraphAfter()
$newCount3 = $d.Paragraphs.Count
# New paragraph 11: new arxiv link
$rng = ($d.Paragraphs($newCount3)).Range
$textRange = $d.Range($rng.Start, $rng.End - 1)
$textRange.Select()
if ($textRange.Start -lt $textRange.End) {
    $word.Selection.Delete()
}
$word.Selection.TypeText('https://arxiv.org/pdf/1803.03635')

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
